$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 6361.8184
$ws.Range("I100").Value = 1868.5714
$ws.Range("J100").Value = 14225
$ws.Range("K100").Value = 1868.5714
$ws.Range("L100").Value = 14225
$ws.Range("M100").Value = -1327.5714
$ws.Range("N100").Value = -15307

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 221097.1
$ws.Range("I32").Value = 221063.27
$ws.Range("J32").Value = 221408.4
$ws.Range("K32").Value = 221063.27
$ws.Range("L32").Value = 221408.4
$ws.Range("M32").Value = -220776.27
$ws.Range("N32").Value = -221982.4

# Row 88: The Mast Chance | Adamantite Rivets
$ws.Range("H88").Value = 19949.916
$ws.Range("I88").Value = 4416.5
$ws.Range("J88").Value = 35483.332
$ws.Range("K88").Value = 4416.5
$ws.Range("L88").Value = 35483.332
$ws.Range("M88").Value = -4010.5
$ws.Range("N88").Value = -36295.332

# Row 91: The Rose and the Riveter (L) | Adamantite Rivets
$ws.Range("H91").Value = 19949.916
$ws.Range("I91").Value = 4416.5
$ws.Range("J91").Value = 35483.332
$ws.Range("K91").Value = 4416.5
$ws.Range("L91").Value = 35483.332
$ws.Range("M91").Value = -3012.5
$ws.Range("N91").Value = -38291.332

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 961.2432
$ws.Range("I97").Value = 722.2963
$ws.Range("J97").Value = 1606.4
$ws.Range("K97").Value = 722.2963
$ws.Range("L97").Value = 1606.4
$ws.Range("M97").Value = -226.2963
$ws.Range("N97").Value = -2598.4

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 439384.06
$ws.Range("I86").Value = 2013.6364
$ws.Range("J86").Value = 1401599
$ws.Range("K86").Value = 2013.6364
$ws.Range("L86").Value = 1401599
$ws.Range("M86").Value = -890.6364000000001
$ws.Range("N86").Value = -1403845

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 439384.06
$ws.Range("I89").Value = 2013.6364
$ws.Range("J89").Value = 1401599
$ws.Range("K89").Value = 10068.182
$ws.Range("L89").Value = 7007995
$ws.Range("M89").Value = -4452.182000000001
$ws.Range("N89").Value = -7019227

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 1041.3572
$ws.Range("I99").Value = 864.9167
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 864.9167
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = 633.0833
$ws.Range("N99").Value = -5096

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 1164.25
$ws.Range("I107").Value = 724.4167
$ws.Range("J107").Value = 2483.75
$ws.Range("K107").Value = 724.4167
$ws.Range("L107").Value = 2483.75
$ws.Range("M107").Value = 1195.5833
$ws.Range("N107").Value = -6323.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 30212.805
$ws.Range("I31").Value = 32101.906
$ws.Range("J31").Value = 15100
$ws.Range("K31").Value = 32101.906
$ws.Range("L31").Value = 15100
$ws.Range("M31").Value = -31806.906
$ws.Range("N31").Value = -15690

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 30212.805
$ws.Range("I34").Value = 32101.906
$ws.Range("J34").Value = 15100
$ws.Range("K34").Value = 32101.906
$ws.Range("L34").Value = 15100
$ws.Range("M34").Value = -31899.906
$ws.Range("N34").Value = -15504

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 593.7
$ws.Range("I107").Value = 422.2
$ws.Range("J107").Value = 765.2
$ws.Range("K107").Value = 422.2
$ws.Range("L107").Value = 765.2
$ws.Range("M107").Value = 1497.8
$ws.Range("N107").Value = -4605.2

$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand | Raisins
$ws.Range("H7").Value = 217.45454
$ws.Range("I7").Value = 177.5
$ws.Range("J7").Value = 240.28572
$ws.Range("K7").Value = 532.5
$ws.Range("L7").Value = 720.85716
$ws.Range("M7").Value = -420.5
$ws.Range("N7").Value = -944.85716

# Row 92: Oh No Udon | Gyr Abanian Flour
$ws.Range("H92").Value = 1281.0834
$ws.Range("I92").Value = 796
$ws.Range("J92").Value = 1523.625
$ws.Range("K92").Value = 2388
$ws.Range("L92").Value = 4570.875
$ws.Range("M92").Value = -1140
$ws.Range("N92").Value = -7066.875

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 569.0741
$ws.Range("I113").Value = 418.57144
$ws.Range("J113").Value = 621.75
$ws.Range("K113").Value = 1255.71432
$ws.Range("L113").Value = 1865.25
$ws.Range("M113").Value = 914.28568
$ws.Range("N113").Value = -6205.25

# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 504.1
$ws.Range("I122").Value = 324.5625
$ws.Range("J122").Value = 1222.25
$ws.Range("K122").Value = 2921.0625
$ws.Range("L122").Value = 11000.25
$ws.Range("M122").Value = -471.0625
$ws.Range("N122").Value = -15900.25

# Row 125: At Any Temperature | Borscht
$ws.Range("H125").Value = 7097.4443
$ws.Range("I125").Value = 2960
$ws.Range("J125").Value = 9166.166999999999
$ws.Range("K125").Value = 8880
$ws.Range("L125").Value = 27498.501
$ws.Range("M125").Value = -3960
$ws.Range("N125").Value = -37338.501

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 4468.524
$ws.Range("I70").Value = 4444.1577
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 4444.1577
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -4174.1577
$ws.Range("N70").Value = -5240

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 4468.524
$ws.Range("I73").Value = 4444.1577
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 4444.1577
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -3508.1577
$ws.Range("N73").Value = -6572

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 3677.0588
$ws.Range("I122").Value = 3457.1304
$ws.Range("J122").Value = 4136.909
$ws.Range("K122").Value = 10371.3912
$ws.Range("L122").Value = 12410.727
$ws.Range("M122").Value = -7921.3912
$ws.Range("N122").Value = -17310.727

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 1943.8572
$ws.Range("I126").Value = 1520
$ws.Range("J126").Value = 2179.3333
$ws.Range("K126").Value = 4560
$ws.Range("L126").Value = 6537.999899999999
$ws.Range("M126").Value = -2090
$ws.Range("N126").Value = -11477.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 38004.57
$ws.Range("I100").Value = 101700.3
$ws.Range("J100").Value = 2618.0557
$ws.Range("K100").Value = 101700.3
$ws.Range("L100").Value = 2618.0557
$ws.Range("M100").Value = -101159.3
$ws.Range("N100").Value = -3700.0557

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 1460.125
$ws.Range("I107").Value = 1156.5
$ws.Range("J107").Value = 1763.75
$ws.Range("K107").Value = 3469.5
$ws.Range("L107").Value = 5291.25
$ws.Range("M107").Value = -1549.5
$ws.Range("N107").Value = -9131.25

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 3542200.5
$ws.Range("I136").Value = 3403338.5
$ws.Range("J136").Value = 5000252.5
$ws.Range("K136").Value = 10210015.5
$ws.Range("L136").Value = 15000757.5
$ws.Range("M136").Value = -10207465.5
$ws.Range("N136").Value = -15005857.5
